# Apply the "update CV1 HGI and fish species. Correct text in priorities doc related to FSW"
# edit to the priorities workbook.
#
# The edit corrects a recurring typo/wording issue ("ander FRPA for" -> "under FRPA
# due to downstream") in several "comments" cells (column H) of the "priorities" sheet,
# and appends a new observation sentence to one of them (H23 - the CV1_us site).
#
# NOTE: cells are edited in this specific order (H6, H19, H29, H31, H23) because the
# underlying shared-strings table is rebuilt (orphaned strings removed/compacted, new
# strings appended in the order they are written) - this order reproduces the exact
# shared-string table layout seen in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("priorities")

# --- H6 : site 125179_us ---------------------------------------------------
$cell = $ws.Cells.Item(6, 8)
$cell.Value = $cell.Text.Replace("ander FRPA for", "under FRPA due to downstream")

# --- H19 : site (2014 FSW-TAG f-7-022, "Some deep pools and boulders") -----
$cell = $ws.Cells.Item(19, 8)
$cell.Value = $cell.Text.Replace( `
    "ander FRPA for Bull Trout and Arctic Grayling (Beaudry", `
    "under FRPA due to downstream Bull Trout and Arctic Grayling  (Beaudry")

# --- H29 : site 125180_us ---------------------------------------------------
$cell = $ws.Cells.Item(29, 8)
$cell.Value = $cell.Text.Replace( `
    "ander FRPA for Bull Trout and Arctic Grayling (Beaudry", `
    "under FRPA due to downstream Bull Trout and Arctic Grayling  (Beaudry")

# --- H31 : site 125186_us ---------------------------------------------------
$cell = $ws.Cells.Item(31, 8)
$cell.Value = $cell.Text.Replace("ander FRPA for", "under FRPA due to downstream")

# --- H23 : site CV1_us -------------------------------------------------------
# correct the wording AND append the new Bull Trout redd observation sentence
$cell = $ws.Cells.Item(23, 8)
$newText = $cell.Text.Replace( `
    "ander FRPA for Bull Trout and Arctic Grayling (Beaudry", `
    "under FRPA due to downstream Bull Trout and Arctic Grayling  (Beaudry")
$newText = $newText + " 200mm Bull Trout (suspected) observed upstream near redd."
$cell.Value = $newText

# --- Column widths: split former shared 7:8 width, widen column H (comments) ---
$ws.Columns.Item(8).ColumnWidth = 87.8

# --- Sheet view: scroll down and select H25, matching the saved view state ----
$ws.Range("H25").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
